$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (shifts existing row 14 down to row 15)
$ws.Rows.Item(14).Insert()

# Fill the new row 14 with the new coordinate pair
$ws.Range("A14").Value = 47.601405800000002
$ws.Range("B14").Value = -122.3235825

# Update the selection to match the saved state
$ws.Range("D20").Select()
